$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row above the current row 218, pushing the existing
# rows 218..325 down to 219..326 (dimension grows from A1:T325 to A1:T326).
$ws.Rows.Item(218).Insert()

# Populate the newly inserted row 218 with the new data record.
$ws.Range("A218").Value = 10
$ws.Range("B218").Value = "Vega Modelo de Temuco"
$ws.Range("C218").Value = "La Araucanía"
$ws.Range("D218").Value = 45205
$ws.Range("E218").Value = 9
$ws.Range("F218").Value = "Fruta"
$ws.Range("G218").Value = 100101
$ws.Range("H218").Value = "Berries"
$ws.Range("I218").Value = 100112025
$ws.Range("J218").Value = "Frutilla"
$ws.Range("K218").Value = "Sin especificar"
$ws.Range("L218").Value = "Primera"
$ws.Range("M218").Value = 650
$ws.Range("N218").Value = 14000
$ws.Range("O218").Value = 15000
$ws.Range("P218").Value = 14615
$ws.Range("Q218").Value = "`$/bandeja 7 kilos"
$ws.Range("R218").Value = "Provincia de Melipilla"
$ws.Range("S218").Value = 2088
$ws.Range("T218").Value = 7
